$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.768714904785156
$ws.Range("B1").Value = 6.19842529296875
$ws.Range("C1").Value = 5.481249332427979
$ws.Range("D1").Value = 6.364148616790771
$ws.Range("E1").Value = 3.731287956237793
